# Update countries & provincias Spain
# Applies the data refresh described by the commit: updates the "Pais" sheet
# with newer COVID figures, swaps two pairs of countries that changed rank
# order (Libia/Hong Kong and Groenlandia/Islas Malvinas), and bumps the
# "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country labels that swapped rows (re-sorted by total cases) ---------
# Row 116 now holds Hong Kong's (updated) figures, row 117 now holds Libia's.
$ws.Range("A116").Value = "Hong Kong"
$ws.Range("A117").Value = "Libia"

# Islas Malvinas / Groenlandia also trade places (identical totals, so this
# only changes which label sits on which row).
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# --- Updated timestamp footer ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 12:20"

# --- Updated numeric figures (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4371992
$ws.Range("C4").Value = 153
$ws.Range("D4").Value = 2090231
$ws.Range("E4").Value = 2131909
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 149852

# Row 34: Oman
$ws.Range("B34").Value = 77058
$ws.Range("C34").Value = 1053
$ws.Range("D34").Value = 57028
$ws.Range("E34").Value = 19637
$ws.Range("G34").Value = 9
$ws.Range("H34").Value = 393

# Row 37: Belgica
$ws.Range("D37").Value = 17439
$ws.Range("E37").Value = 38766

# Row 64: Uzbekistan
$ws.Range("B64").Value = 20820
$ws.Range("C64").Value = 289
$ws.Range("E64").Value = 9599

# Row 65: Austria
$ws.Range("B65").Value = 20558
$ws.Range("C65").Value = 86
$ws.Range("D65").Value = 18246
$ws.Range("E65").Value = 1599
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 713

# Row 84: Senegal
$ws.Range("B84").Value = 9764
$ws.Range("C84").Value = 83
$ws.Range("D84").Value = 6477
$ws.Range("E84").Value = 3093
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 194

# Row 87: Malasia
$ws.Range("B87").Value = 8904
$ws.Range("C87").Value = 7
$ws.Range("D87").Value = 8601
$ws.Range("E87").Value = 179

# Row 89: Finlandia
$ws.Range("B89").Value = 7398
$ws.Range("C89").Value = 5
$ws.Range("E89").Value = 149

# Row 116: now Hong Kong
$ws.Range("B116").Value = 2779
$ws.Range("C116").Value = 145
$ws.Range("D116").Value = 1495
$ws.Range("E116").Value = 1264
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 20

# Row 117: now Libia
$ws.Range("B117").Value = 2669
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 553
$ws.Range("E117").Value = 2056
$ws.Range("G117").Value = 18
$ws.Range("H117").Value = 60

# Row 125: Eslovenia
$ws.Range("B125").Value = 2087
$ws.Range("C125").Value = 5
$ws.Range("D125").Value = 1733
$ws.Range("E125").Value = 238
